# Generate Report for Handback
# Updates the "Ready for handoff" rows (for the e2e/9f27646a... file) to reflect
# that the handback has now completed ("Handed back: in sync with en-US"),
# records the actual handback timestamps, and clears the stale
# "version not latest" error message now that the handback is current.

$wb = $excel.ActiveWorkbook

$statusDone = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E3").Value = $statusDone
$ovw.Range("F3").Value = $statusDone

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusDone
$zhcn.Range("K3").Value = "2016-08-24 10:49:52"
$zhcn.Range("P3").Value = ""

# --- de-de sheet ------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusDone
$dede.Range("K3").Value = "2016-08-24 10:49:59"
$dede.Range("P3").Value = ""

# The Error Detail column no longer needs to be wide enough for a long
# error message now that it's empty for this row.
$zhcn.Columns.Item(16).ColumnWidth = 13.7470528738839
$dede.Columns.Item(16).ColumnWidth = 13.7470528738839
